$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers N1:P1 (same style as other header cells, reuse M1's formatting)
$ws.Range("N1").Value = "M_Vitorias"
$ws.Range("O1").Value = "M_Derrotas"
$ws.Range("P1").Value = "M_Empates"
$ws.Range("M1").Copy()
$ws.Range("N1:P1").PasteSpecial(-4122)

# Fill in computed values for rows 2 through 50 (averages per "baba"):
# N = Total_Vitorias (D) / Total_Babas (H)
# O = Total_Derrotas (E) / Total_Babas (H)
# P = Total_Empates (F) / Total_Babas (H)
for ($r = 2; $r -le 50; $r++) {
    $d = $ws.Cells.Item($r, 4).Value2
    $e = $ws.Cells.Item($r, 5).Value2
    $f = $ws.Cells.Item($r, 6).Value2
    $h = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 14).Value = $d / $h
    $ws.Cells.Item($r, 15).Value = $e / $h
    $ws.Cells.Item($r, 16).Value = $f / $h
}
